$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Chris K." (2nd sheet): split the combined "Premium Cat Food" /
# "Brush" out-of-stock note into its own proper line item row, same as the
# pattern already used on the "Beth S." sheet.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Row 2: remove the "Out of Stock" note (E2) and update the Subtotal/Total
# now that "Premium Cat Food" is its own purchased line below.
$ws2.Range("E2").ClearContents()
$ws2.Range("I2").Value = 33.97
$ws2.Range("K2").Value = 54.17

# Insert a fresh row at 4 (pushes the old row 4 "Oatmeal Soap" down to row 5)
# so we can split what used to be row 3 into two rows.
$ws2.Rows.Item(4).Insert()

# Row 3 becomes the "Premium Cat Food" purchase (previously the Out of Stock
# note in E2), with its real cost.
$ws2.Range("A3").Value = "Premium Cat Food"
$ws2.Range("C3").Value = 11.99

# New row 4 keeps the original row-3 "Brush" line item.
$ws2.Range("A4").Value = "Brush"
$ws2.Range("B4").Value = "Pet"
$ws2.Range("C4").Value = 4.99

# ---------------------------------------------------------------------------
# Sheet "Mary M." (3rd sheet): same treatment - split out the "Fruits" /
# "Vegetables" notes that were crammed into column E into their own rows.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("E2").ClearContents()
$ws3.Range("I2").Value = 32.979999999999997
$ws3.Range("K2").Value = 52.98

# Row 3 becomes the "Fruits" purchase (previously the Out of Stock note in
# E2), replacing the old "Basic Dog Food" line which moves further down.
# Also drop the old "Vegetables" note that lived in E3 - it becomes its own
# row below instead.
$ws3.Range("A3").Value = "Fruits"
$ws3.Range("C3").Value = 7
$ws3.Range("E3").ClearContents()

# Insert two fresh rows at 4 and 5, pushing nothing further down (sheet ends
# at row 3 before this edit).
$ws3.Range("A4").Value = "Vegetables"
$ws3.Range("B4").Value = "Grocery"
$ws3.Range("C4").Value = 5

$ws3.Range("A5").Value = "Basic Dog Food"
$ws3.Range("B5").Value = "Grocery"
$ws3.Range("C5").Value = 8.99
